# Generate Report for handback
#
# This script mirrors a localization "handback" run: the files that were
# handed off (a.md.md / b.md.md, status "Ready for handoff") have now come
# back from the translators in sync with en-US, so:
#   - the status text changes everywhere it appears (Overview + per-locale
#     sheets) from "Ready for handoff" to "Handed back: in sync with en-US"
#   - each per-locale sheet grows two new populated columns for the rows
#     that were actually included in the handoff: "Latest Target File" (E)
#     and "Latest Handback File" (F), each a hyperlink to the relevant file
#   - "Latest Handback DateTime" (G) is stamped with the real handback time
#     instead of the epoch placeholder "0001-01-01 00:00:00"

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: just the status label changes (same two localized
# columns get the new text for each in-scope file row).
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# Per-locale sheets: zh-cn and de-de get the same treatment.
# ---------------------------------------------------------------------
$locales = @(
    @{
        Sheet = "zh-cn"
        XlfName = "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf"
        HandbackDateTime = "2016-02-16 09:47:03"
        TargetUrl = "https://github.com/OpenLocalizationTest/oltest/blob/170e2003dbbf9ddea08c3e79e0dc84b95449209b/e2e/a.md.md"
        HandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a782058e2f9dfd19f5fdc712bb334dfcdb85b578/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/hb/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf"
    },
    @{
        Sheet = "de-de"
        XlfName = "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf"
        HandbackDateTime = "2016-02-16 09:47:29"
        TargetUrl = "https://github.com/OpenLocalizationTest/oltest/blob/170e2003dbbf9ddea08c3e79e0dc84b95449209b/e2e/a.md.md"
        HandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/28cab5b4d2114b2d9b254726d701b708646c1d73/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/hb/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf"
    }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    # Rows 2 and 3 are the two files ("a.md.md"/"b.md.md") that were
    # included in the handoff (Handoff Reason = "Include"); row 4 is the
    # ".localization-config" file, which is never localized and is left
    # untouched.
    foreach ($row in 2, 3) {
        $ws.Cells.Item($row, 2).Value = $newStatus   # B: Status

        # E: Latest Target File - the localized file that came back
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 5), $locale.TargetUrl, "", "", "a.md.md")

        # F: Latest Handback File - the handback .xlf for this locale
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $locale.HandbackUrl, "", "", $locale.XlfName)

        # G: Latest Handback DateTime - stamp the real handback time
        $ws.Cells.Item($row, 7).Value = $locale.HandbackDateTime
    }
}
